# Add 2022-Q3 data
#
# Starting state:
#   Sheet1 "总计"    - summary sheet, row2 holds the latest quarter's totals
#   Sheet2 "2022-Q2" - per-fund detail for 2022-Q2
#
# Target state:
#   Sheet1 "总计"    - row2 now holds 2022-Q3 totals, row3 (new) holds the
#                      2022-Q2 totals that used to live in row2
#   Sheet2 "2022-Q3" - NEW sheet (inserted before "2022-Q2") with the
#                      per-fund detail for 2022-Q3
#   Sheet3 "2022-Q2" - the original detail sheet, content untouched, just
#                      pushed one slot to the right

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. "总计" summary sheet: push the existing 2022-Q2 row down to row 3,
#    then overwrite row 2 with the new 2022-Q3 totals.
# ---------------------------------------------------------------------
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)   # xlPasteFormats - carry the A-column style down

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 0.1

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.09

# ---------------------------------------------------------------------
# 2. New "2022-Q3" detail sheet, inserted right before "2022-Q2".
#    Built by duplicating "总计" (so it inherits the same page setup /
#    header style) and then replacing its contents.
# ---------------------------------------------------------------------
$summary.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"
$q3Sheet.Cells.Clear()

# Re-apply the header / index-column formatting that "总计" used.
$summary.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$summary.Range("A2").Copy()
$q3Sheet.Range("A2:A3").PasteSpecial(-4122)

$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").NumberFormat = "@"
$q3Sheet.Range("B2").Value = "501021"
$q3Sheet.Range("B2").Style = "Normal"
$q3Sheet.Range("C2").Value = "华宝标普香港上市中国中小盘指数（LOF）A"
$q3Sheet.Range("D2").NumberFormat = "@"
$q3Sheet.Range("D2").Value = "4.19"
$q3Sheet.Range("D2").Style = "Normal"
$q3Sheet.Range("E2").NumberFormat = "@"
$q3Sheet.Range("E2").Value = "92.99"
$q3Sheet.Range("E2").Style = "Normal"
$q3Sheet.Range("F2").NumberFormat = "@"
$q3Sheet.Range("F2").Value = "2.06"
$q3Sheet.Range("F2").Style = "Normal"
$q3Sheet.Range("G2").NumberFormat = "@"
$q3Sheet.Range("G2").Value = "0.0863"
$q3Sheet.Range("G2").Style = "Normal"
$q3Sheet.Range("H2").Value = 5

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").NumberFormat = "@"
$q3Sheet.Range("B3").Value = "006127"
$q3Sheet.Range("B3").Style = "Normal"
$q3Sheet.Range("C3").Value = "华宝标普香港上市中国中小盘指数（LOF）C"
$q3Sheet.Range("D3").NumberFormat = "@"
$q3Sheet.Range("D3").Value = "0.24"
$q3Sheet.Range("D3").Style = "Normal"
$q3Sheet.Range("E3").NumberFormat = "@"
$q3Sheet.Range("E3").Value = "92.99"
$q3Sheet.Range("E3").Style = "Normal"
$q3Sheet.Range("F3").NumberFormat = "@"
$q3Sheet.Range("F3").Value = "2.06"
$q3Sheet.Range("F3").Style = "Normal"
$q3Sheet.Range("G3").NumberFormat = "@"
$q3Sheet.Range("G3").Value = "0.0049"
$q3Sheet.Range("G3").Style = "Normal"
$q3Sheet.Range("H3").Value = 5
